$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Secret_Child_Name (C) and Secret_Child_EmailID (D) columns
# for rows 2-16 to reflect the corrected Secret Santa assignments.

$ws.Cells.Item(2, 3).Value = "Matthew King"
$ws.Cells.Item(2, 4).Value = "matthew.king.jr@acme.com"

$ws.Cells.Item(3, 3).Value = "Isabella Scott"
$ws.Cells.Item(3, 4).Value = "isabella.scott@acme.com"

$ws.Cells.Item(4, 3).Value = "Charlie Ross"
$ws.Cells.Item(4, 4).Value = "charlie.ross.jr@acme.com"

$ws.Cells.Item(5, 3).Value = "Charlie Ross"
$ws.Cells.Item(5, 4).Value = "charlie.ross@acme.com"

$ws.Cells.Item(6, 3).Value = "Piper Stewart"
$ws.Cells.Item(6, 4).Value = "piper.stewart@acme.com"

$ws.Cells.Item(7, 3).Value = "Charlie Wright"
$ws.Cells.Item(7, 4).Value = "charlie.wright@acme.com"

$ws.Cells.Item(8, 3).Value = "Matthew King"
$ws.Cells.Item(8, 4).Value = "matthew.king@acme.com"

$ws.Cells.Item(9, 3).Value = "Benjamin Collins"
$ws.Cells.Item(9, 4).Value = "benjamin.collins@acme.com"

$ws.Cells.Item(10, 3).Value = "Mark Lawrence"
$ws.Cells.Item(10, 4).Value = "mark.lawrence@acme.com"

$ws.Cells.Item(11, 3).Value = "Hamish Murray"
$ws.Cells.Item(11, 4).Value = "hamish.murray.jr@acme.com"

$ws.Cells.Item(12, 3).Value = "Spencer Allen"
$ws.Cells.Item(12, 4).Value = "spencer.allen@acme.com"

$ws.Cells.Item(13, 3).Value = "Hamish Murray"
$ws.Cells.Item(13, 4).Value = "hamish.murray@acme.com"

$ws.Cells.Item(14, 3).Value = "Hamish Murray"
$ws.Cells.Item(14, 4).Value = "hamish.murray.sr@acme.com"

$ws.Cells.Item(15, 3).Value = "Ethan Murray"
$ws.Cells.Item(15, 4).Value = "ethan.murray@acme.com"

$ws.Cells.Item(16, 3).Value = "Layla Graham"
$ws.Cells.Item(16, 4).Value = "layla.graham@acme.com"
